$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date serial for every
# data row (rows 2-391). The whole column was bumped by one day
# (46074 -> 46075), i.e. from 2026-02-21 to 2026-02-22.
$ws.Range("C2:C391").Value = 46075
